$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Title ---
Replace-Text "Unveiling the Enigma of Dark Matter" "Government: The Balancing Act of Society"

# --- Byline (name) ---
Replace-Text "Dr. Alan Cassidy" "Emma Anderson"

# --- Email address ---
Replace-Text "alan" "emma"
Replace-Text "cassidy09@gmail" "anderson@validmail"

# --- Body paragraph 1 ---
Replace-Text "The universe, a boundless tapestry of celestial wonders, conceals an enigmatic entity known as dark matter, an invisible yet pervasive force that wields gravitational influence without emitting any light." "In the tapestry of human existence, the thread of government runs through civilizations like a leitmotif, binding people together in a symphony of laws, rights, and responsibilities."

Replace-Text " Comprising nearly 85% of the universe's total mass, it remains an elusive puzzle, tantalizing and confounding scientists worldwide." " It sculpts the landscape of our societies, determining how we interact with each other, and shaping our collective destiny. Government is the conductor of our communal orchestra, harmonizing the diverse notes of individual aspirations into a resonant melody of societal progress. Its architects are the composers who craft the score of our laws, the conductors who guide our institutions, and the players who bring life to the policies that shape our lives."

Replace-Text "Through gravitational lensing and its impact on the rotation of galaxies, the existence of dark matter has been inferred." "The story of government is as old as civilization itself."

Replace-Text " Galaxies, like celestial whirlpools, spin with an unexpected velocity, exceeding the speed that would be anticipated based on their visible mass alone." " From the earliest tribal councils to the modern-day nation-states, the need for governance has been a constant, a testament to our innate desire for order and cooperation."

Replace-Text " This discrepancy suggests the presence of unseen matter, exerting a gravitational pull that governs the galaxy's rotation." " Over time, governments have evolved in myriad forms, each reflecting the unique circumstances and aspirations of the people they serve. From monarchies to democracies, from totalitarian regimes to constitutional republics, the tapestry of government has been woven with both triumphs and failures, successes and tribulations."

Replace-Text "Furthermore, dark matter's influence is evident in the behavior of galaxy clusters, vast congregations of galaxies bound together by gravity." "Government is like a living organism, constantly adapting and evolving as it responds to the ever-changing landscape of society."

Replace-Text " The motion of galaxies within these clusters defies conventional expectations, indicating the presence of significantly more mass than what is visible." " It is both a mirror, reflecting the values and aspirations of its people, and a catalyst for change, driving progress and innovation."

Replace-Text " This unseen mass, invisible to telescopes and instruments, has thus far remained shrouded in mystery." " The story of government is a story of people, of their struggles and triumphs, their hopes and fears, their dreams and aspirations. It is a story that is still being written, a story that we, as citizens, have the privilege and responsibility to shape."

# --- Summary paragraph ---
Replace-Text "The existence of dark matter, a mysterious and unseen entity, is supported by various lines of evidence, including gravitational lensing, the velocity of galaxies, and the behavior of galaxy clusters." "This essay explored the multifaceted nature of government, its historical evolution, and its profound impact on society."

Replace-Text " Despite its pervasive presence, dark matter remains elusive and enigmatic, challenging our understanding of the universe." " It delved into the various forms of government, from monarchies to democracies, highlighting the unique characteristics and challenges of each."

Replace-Text " Its composition, properties, and role in the grand cosmic scheme continue to puzzle scientists, beckoning them to unravel the secrets of this hidden realm." " The essay emphasized the importance of citizen participation and engagement in governance, as well as the need for governments to be responsive and accountable to the people they serve. Ultimately, it underscored the enduring significance of government as the foundation of order, cooperation, and progress in human society."

# --- Trailing empty paragraph ---
$d.Content.InsertParagraphAfter() | Out-Null
